$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-25 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-26 Friday", 2) | Out-Null
$d.Content.Find.Execute("86-83=3", $true, $false, $false, $false, $false, $true, 1, $false, "58+12=70", 2) | Out-Null
$d.Content.Find.Execute("97-38=59", $true, $false, $false, $false, $false, $true, 1, $false, "19+48=67", 2) | Out-Null
$d.Content.Find.Execute("77-66=11", $true, $false, $false, $false, $false, $true, 1, $false, "45-37=8", 2) | Out-Null
$d.Content.Find.Execute("59-21=38", $true, $false, $false, $false, $false, $true, 1, $false, "62-38=24", 2) | Out-Null
$d.Content.Find.Execute("51+21=72", $true, $false, $false, $false, $false, $true, 1, $false, "5+43=48", 2) | Out-Null
$d.Content.Find.Execute("72-9=63", $true, $false, $false, $false, $false, $true, 1, $false, "55+31=86", 2) | Out-Null
$d.Content.Find.Execute("98-19=79", $true, $false, $false, $false, $false, $true, 1, $false, "5+81=86", 2) | Out-Null
$d.Content.Find.Execute("0+80=80", $true, $false, $false, $false, $false, $true, 1, $false, "73-9=64", 2) | Out-Null
$d.Content.Find.Execute("2+76=78", $true, $false, $false, $false, $false, $true, 1, $false, "35+24=59", 2) | Out-Null
$d.Content.Find.Execute("10+60=70", $true, $false, $false, $false, $false, $true, 1, $false, "25-14=11", 2) | Out-Null
$d.Content.Find.Execute("81-33=48", $true, $false, $false, $false, $false, $true, 1, $false, "99-9=90", 2) | Out-Null
$d.Content.Find.Execute("41-12=29", $true, $false, $false, $false, $false, $true, 1, $false, "71+3=74", 2) | Out-Null
$d.Content.Find.Execute("69-23=46", $true, $false, $false, $false, $false, $true, 1, $false, "94-11=83", 2) | Out-Null
$d.Content.Find.Execute("92-75=17", $true, $false, $false, $false, $false, $true, 1, $false, "25+68=93", 2) | Out-Null
$d.Content.Find.Execute("5+8=13", $true, $false, $false, $false, $false, $true, 1, $false, "34+10=44", 2) | Out-Null
$d.Content.Find.Execute("10+21=31", $true, $false, $false, $false, $false, $true, 1, $false, "73-17=56", 2) | Out-Null
$d.Content.Find.Execute("70-52=18", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=65", 2) | Out-Null
$d.Content.Find.Execute("39+10=49", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=41", 2) | Out-Null
$d.Content.Find.Execute("56-29=27", $true, $false, $false, $false, $false, $true, 1, $false, "82-39=43", 2) | Out-Null
$d.Content.Find.Execute("65+25=90", $true, $false, $false, $false, $false, $true, 1, $false, "82-82=0", 2) | Out-Null
$d.Content.Find.Execute("81-44=37", $true, $false, $false, $false, $false, $true, 1, $false, "75-53=22", 2) | Out-Null
$d.Content.Find.Execute("44-2=42", $true, $false, $false, $false, $false, $true, 1, $false, "47+17=64", 2) | Out-Null
$d.Content.Find.Execute("60-0=60", $true, $false, $false, $false, $false, $true, 1, $false, "33+53=86", 2) | Out-Null
$d.Content.Find.Execute("81-34=47", $true, $false, $false, $false, $false, $true, 1, $false, "30+52=82", 2) | Out-Null
$d.Content.Find.Execute("2+3=5", $true, $false, $false, $false, $false, $true, 1, $false, "59+36=95", 2) | Out-Null
$d.Content.Find.Execute("98-50=48", $true, $false, $false, $false, $false, $true, 1, $false, "70-9=61", 2) | Out-Null
$d.Content.Find.Execute("57-53=4", $true, $false, $false, $false, $false, $true, 1, $false, "53+45=98", 2) | Out-Null
$d.Content.Find.Execute("94-88=6", $true, $false, $false, $false, $false, $true, 1, $false, "76-58=18", 2) | Out-Null
$d.Content.Find.Execute("46-21=25", $true, $false, $false, $false, $false, $true, 1, $false, "55-38=17", 2) | Out-Null
$d.Content.Find.Execute("97-32=65", $true, $false, $false, $false, $false, $true, 1, $false, "96-76=20", 2) | Out-Null
$d.Content.Find.Execute("36-26=10", $true, $false, $false, $false, $false, $true, 1, $false, "32-25=7", 2) | Out-Null
$d.Content.Find.Execute("31+63=94", $true, $false, $false, $false, $false, $true, 1, $false, "8+91=99", 2) | Out-Null
$d.Content.Find.Execute("95+0=95", $true, $false, $false, $false, $false, $true, 1, $false, "62+29=91", 2) | Out-Null
$d.Content.Find.Execute("85-65=20", $true, $false, $false, $false, $false, $true, 1, $false, "82-13=69", 2) | Out-Null
$d.Content.Find.Execute("80-27=53", $true, $false, $false, $false, $false, $true, 1, $false, "31+60=91", 2) | Out-Null
$d.Content.Find.Execute("89-76=13", $true, $false, $false, $false, $false, $true, 1, $false, "30+13=43", 2) | Out-Null
$d.Content.Find.Execute("30+31=61", $true, $false, $false, $false, $false, $true, 1, $false, "99-64=35", 2) | Out-Null
$d.Content.Find.Execute("82-64=18", $true, $false, $false, $false, $false, $true, 1, $false, "56-2=54", 2) | Out-Null
$d.Content.Find.Execute("50+44=94", $true, $false, $false, $false, $false, $true, 1, $false, "13+86=99", 2) | Out-Null
$d.Content.Find.Execute("8+73=81", $true, $false, $false, $false, $false, $true, 1, $false, "65-3=62", 2) | Out-Null
$d.Content.Find.Execute("50+43=93", $true, $false, $false, $false, $false, $true, 1, $false, "53-14=39", 2) | Out-Null
$d.Content.Find.Execute("17+21=38", $true, $false, $false, $false, $false, $true, 1, $false, "78-40=38", 2) | Out-Null
$d.Content.Find.Execute("91-78=13", $true, $false, $false, $false, $false, $true, 1, $false, "90-10=80", 2) | Out-Null
$d.Content.Find.Execute("44-4=40", $true, $false, $false, $false, $false, $true, 1, $false, "12+48=60", 2) | Out-Null
$d.Content.Find.Execute("71-26=45", $true, $false, $false, $false, $false, $true, 1, $false, "29+22=51", 2) | Out-Null
$d.Content.Find.Execute("31-24=7", $true, $false, $false, $false, $false, $true, 1, $false, "28+68=96", 2) | Out-Null
$d.Content.Find.Execute("25+67=92", $true, $false, $false, $false, $false, $true, 1, $false, "55-38=17", 2) | Out-Null
$d.Content.Find.Execute("2+12=14", $true, $false, $false, $false, $false, $true, 1, $false, "47+9=56", 2) | Out-Null
$d.Content.Find.Execute("77-40=37", $true, $false, $false, $false, $false, $true, 1, $false, "84-47=37", 2) | Out-Null
$d.Content.Find.Execute("15+41=56", $true, $false, $false, $false, $false, $true, 1, $false, "6+57=63", 2) | Out-Null
$d.Content.Find.Execute("3+37=40", $true, $false, $false, $false, $false, $true, 1, $false, "6+43=49", 2) | Out-Null
$d.Content.Find.Execute("20-2=18", $true, $false, $false, $false, $false, $true, 1, $false, "24+39=63", 2) | Out-Null
$d.Content.Find.Execute("68-38=30", $true, $false, $false, $false, $false, $true, 1, $false, "0+71=71", 2) | Out-Null
$d.Content.Find.Execute("59-35=24", $true, $false, $false, $false, $false, $true, 1, $false, "19+58=77", 2) | Out-Null
$d.Content.Find.Execute("12+85=97", $true, $false, $false, $false, $false, $true, 1, $false, "99-17=82", 2) | Out-Null
$d.Content.Find.Execute("31+7=38", $true, $false, $false, $false, $false, $true, 1, $false, "67+29=96", 2) | Out-Null
$d.Content.Find.Execute("4+93=97", $true, $false, $false, $false, $false, $true, 1, $false, "26+7=33", 2) | Out-Null
$d.Content.Find.Execute("88-37=51", $true, $false, $false, $false, $false, $true, 1, $false, "9+44=53", 2) | Out-Null
$d.Content.Find.Execute("75-59=16", $true, $false, $false, $false, $false, $true, 1, $false, "1+22=23", 2) | Out-Null
$d.Content.Find.Execute("97-97=0", $true, $false, $false, $false, $false, $true, 1, $false, "77-70=7", 2) | Out-Null
$d.Content.Find.Execute("63-3=60", $true, $false, $false, $false, $false, $true, 1, $false, "22+8=30", 2) | Out-Null
$d.Content.Find.Execute("56+33=89", $true, $false, $false, $false, $false, $true, 1, $false, "33+25=58", 2) | Out-Null
$d.Content.Find.Execute("13-5=8", $true, $false, $false, $false, $false, $true, 1, $false, "11+88=99", 2) | Out-Null
$d.Content.Find.Execute("29+55=84", $true, $false, $false, $false, $false, $true, 1, $false, "59-53=6", 2) | Out-Null
$d.Content.Find.Execute("53-39=14", $true, $false, $false, $false, $false, $true, 1, $false, "66-25=41", 2) | Out-Null
$d.Content.Find.Execute("84-60=24", $true, $false, $false, $false, $false, $true, 1, $false, "72-41=31", 2) | Out-Null
$d.Content.Find.Execute("1+60=61", $true, $false, $false, $false, $false, $true, 1, $false, "69-59=10", 2) | Out-Null
$d.Content.Find.Execute("63-12=51", $true, $false, $false, $false, $false, $true, 1, $false, "54+41=95", 2) | Out-Null
$d.Content.Find.Execute("58-13=45", $true, $false, $false, $false, $false, $true, 1, $false, "22+71=93", 2) | Out-Null
$d.Content.Find.Execute("83+13=96", $true, $false, $false, $false, $false, $true, 1, $false, "1+27=28", 2) | Out-Null
$d.Content.Find.Execute("48-38=10", $true, $false, $false, $false, $false, $true, 1, $false, "48-6=42", 2) | Out-Null
$d.Content.Find.Execute("19+23=42", $true, $false, $false, $false, $false, $true, 1, $false, "64-55=9", 2) | Out-Null
$d.Content.Find.Execute("3+33=36", $true, $false, $false, $false, $false, $true, 1, $false, "48+45=93", 2) | Out-Null
$d.Content.Find.Execute("34+17=51", $true, $false, $false, $false, $false, $true, 1, $false, "13+85=98", 2) | Out-Null
$d.Content.Find.Execute("10+77=87", $true, $false, $false, $false, $false, $true, 1, $false, "7+52=59", 2) | Out-Null
$d.Content.Find.Execute("45+19=64", $true, $false, $false, $false, $false, $true, 1, $false, "34-4=30", 2) | Out-Null
$d.Content.Find.Execute("86+3=89", $true, $false, $false, $false, $false, $true, 1, $false, "65-35=30", 2) | Out-Null
$d.Content.Find.Execute("50-3=47", $true, $false, $false, $false, $false, $true, 1, $false, "93-23=70", 2) | Out-Null
$d.Content.Find.Execute("94-34=60", $true, $false, $false, $false, $false, $true, 1, $false, "31-13=18", 2) | Out-Null
$d.Content.Find.Execute("86-46=40", $true, $false, $false, $false, $false, $true, 1, $false, "27-10=17", 2) | Out-Null
$d.Content.Find.Execute("83-24=59", $true, $false, $false, $false, $false, $true, 1, $false, "72-34=38", 2) | Out-Null
$d.Content.Find.Execute("45+7=52", $true, $false, $false, $false, $false, $true, 1, $false, "17+33=50", 2) | Out-Null
$d.Content.Find.Execute("93-91=2", $true, $false, $false, $false, $false, $true, 1, $false, "17+53=70", 2) | Out-Null
$d.Content.Find.Execute("61+4=65", $true, $false, $false, $false, $false, $true, 1, $false, "74-44=30", 2) | Out-Null
$d.Content.Find.Execute("96-86=10", $true, $false, $false, $false, $false, $true, 1, $false, "39-18=21", 2) | Out-Null
$d.Content.Find.Execute("61-26=35", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2) | Out-Null
$d.Content.Find.Execute("73+13=86", $true, $false, $false, $false, $false, $true, 1, $false, "93-36=57", 2) | Out-Null
$d.Content.Find.Execute("53+12=65", $true, $false, $false, $false, $false, $true, 1, $false, "94-64=30", 2) | Out-Null
$d.Content.Find.Execute("21+68=89", $true, $false, $false, $false, $false, $true, 1, $false, "84-4=80", 2) | Out-Null
$d.Content.Find.Execute("38-37=1", $true, $false, $false, $false, $false, $true, 1, $false, "23+21=44", 2) | Out-Null
$d.Content.Find.Execute("3+19=22", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=52", 2) | Out-Null
$d.Content.Find.Execute("69-28=41", $true, $false, $false, $false, $false, $true, 1, $false, "96-25=71", 2) | Out-Null
$d.Content.Find.Execute("68-3=65", $true, $false, $false, $false, $false, $true, 1, $false, "12-11=1", 2) | Out-Null
$d.Content.Find.Execute("78-52=26", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("76+10=86", $true, $false, $false, $false, $false, $true, 1, $false, "14+27=41", 2) | Out-Null
$d.Content.Find.Execute("72-42=30", $true, $false, $false, $false, $false, $true, 1, $false, "5+39=44", 2) | Out-Null
$d.Content.Find.Execute("0+31=31", $true, $false, $false, $false, $false, $true, 1, $false, "84+3=87", 2) | Out-Null
$d.Content.Find.Execute("84-63=21", $true, $false, $false, $false, $false, $true, 1, $false, "98-31=67", 2) | Out-Null
$d.Content.Find.Execute("45-21=24", $true, $false, $false, $false, $false, $true, 1, $false, "48+21=69", 2) | Out-Null
$d.Content.Find.Execute("68-18=50", $true, $false, $false, $false, $false, $true, 1, $false, "68-56=12", 2) | Out-Null
